$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.335.02"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.936.54"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7520"
$ws.Range("E5").Value = "  +5.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.30"
$ws.Range("E6").Value = "  -2.31%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +1.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3196"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07046"
$ws.Range("E10").Value = "  -0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7841"
$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08044"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.931.33"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.416"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.27"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.50"
$ws.Range("E16").Value = "  -1.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.334.72"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.087"
$ws.Range("E18").Value = "  +5.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.93"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008029"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.181.46"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.715"
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.567"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.97"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.12"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1306"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.218"
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.535"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.461"
$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.155"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.342"
$ws.Range("E34").Value = "  +5.83%  "

$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7591"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.783"
$ws.Range("E37").Value = "  +0.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.809"
$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.19"
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.566"
$ws.Range("E41").Value = "  +3.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4528"
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.987"
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8373"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.725"
$ws.Range("E46").Value = "  +3.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.03"
$ws.Range("E47").Value = "  +3.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.75"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.83"
$ws.Range("E49").Value = "  +4.00%  "

$ws.Range("E50").Value = "  +8.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "972.35"
$ws.Range("E51").Value = "  +5.91%  "
